$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ZNbVv125"
$ws.Range("B2").Value = 23090735
$ws.Range("C2").Value = "atucvvd74"
$ws.Range("D2").Value = "u4#QqU!7"
$ws.Range("F2").Value = "yMQASdYk"
$ws.Range("G2").Value = "TwYw"
